# "Add case 4 and 5 / Push before adding seeds to struct"
#
# Updates the existing "Case 4 - study 1 (pc)" block (rows 15-30) with a few
# refreshed values, then adds a brand new "Case 4 - study 2" block in rows
# 32-49 with two side-by-side variants: "optimal" (cols A:B) and "test"
# (cols G:H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = (Get-Date -Year 2022 -Month 7 -Day 14).Date

# ---------------------------------------------------------------------------
# Existing "Case 4 - study 1 (pc)" block: a handful of values were refreshed.
# ---------------------------------------------------------------------------
$ws.Range("B18").Value = "10 m^3/s"
$ws.Range("B23").Value = $newDate
$ws.Range("B26").Value = 0.2
$ws.Range("B27").Value = 100
$ws.Range("B29").Value = 0.3
$ws.Range("B30").Value = "10 md"

# ---------------------------------------------------------------------------
# New "Case 4 - study 2" block (row 31 stays blank, just like the gap before
# the study-1 block at row 13).
# ---------------------------------------------------------------------------
$ws.Range("A32").Value = "Case 4 - study 2 (optimal)"
$ws.Range("G32").Value = "Case 4 - study 2 (test)"

$ws.Range("A33").Value = "nr layers"
$ws.Range("B33").Value = "15, 30, 45"
$ws.Range("G33").Value = "nr layers"
$ws.Range("H33").Value = 30

$ws.Range("A34").Value = "[nx, nz]"
$ws.Range("B34").Value = "[75, 60]"
$ws.Range("G34").Value = "[nx, nz]"
$ws.Range("H34").Value = "[50, 40]"

$ws.Range("A35").Value = "[lx, lz]"
$ws.Range("B35").Value = "[800, 400]"
$ws.Range("G35").Value = "[lx, lz]"
$ws.Range("H35").Value = "[800, 400]"

$ws.Range("A36").Value = "well_rate"
$ws.Range("B36").Value = "sum(poreVolume)/100 / (inj_stop*tot_time)"
$ws.Range("G36").Value = "well_rate"
$ws.Range("H36").Value = "sum(poreVolume)/100 / (inj_stop*tot_time)"

$ws.Range("A37").Value = "inj_stop"
$ws.Range("B37").Value = 0.3
$ws.Range("G37").Value = "inj_stop"
$ws.Range("H37").Value = 0.2

$ws.Range("A38").Value = "rampup"
$ws.Range("B38").Value = "(16000 days, 75 days, 10)"
$ws.Range("G38").Value = "rampup"
$ws.Range("H38").Value = "(500 years, 200 days, 10)"

$ws.Range("A39").Value = "init_pres"
$ws.Range("B39").Value = "100*barsa"
$ws.Range("G39").Value = "init_pres"
$ws.Range("H39").Value = "100*barsa"

$ws.Range("A40").Value = "{layer: seeds}"
$ws.Range("B40").Value = "{15: [6551, 1167], 30: [5223] }"
$ws.Range("G40").Value = "{layer: seeds}"
$ws.Range("H40").Value = "{30: [3042] }"

# Copy the existing date-styled cell (B23) onto B41/H41 first so they pick up
# the same number format (style index), then overwrite with the real value -
# otherwise a fresh .Value = <date> assignment invents a brand-new style.
$ws.Range("A41").Value = "date"
$ws.Range("B23").Copy($ws.Range("B41"))
$ws.Range("B41").Value = $newDate
$ws.Range("G41").Value = "date"
$ws.Range("B23").Copy($ws.Range("H41"))
$ws.Range("H41").Value = $newDate

$ws.Range("A42").Value = "p_e"
$ws.Range("B42").Value = "0.5*barsa"
$ws.Range("G42").Value = "pc_median"
$ws.Range("H42").Value = "1*barsa"

$ws.Range("A43").Value = "p_cap"
$ws.Range("B43").Value = "3*barsa"
$ws.Range("G43").Value = "std_gauss"
$ws.Range("H43").Value = 0.8

$ws.Range("A44").Value = "std_gauss"
$ws.Range("B44").Value = 0.8
$ws.Range("G44").Value = "corr_len_x"
$ws.Range("H44").Value = 200

$ws.Range("A45").Value = "corr_len_x"
$ws.Range("B45").Value = 200
$ws.Range("G45").Value = "corr_len_z"
$ws.Range("H45").Value = 10

$ws.Range("A46").Value = "corr_len_z"
$ws.Range("B46").Value = 10
$ws.Range("G46").Value = "poro"
$ws.Range("H46").Value = 0.3

$ws.Range("A47").Value = "poro"
$ws.Range("B47").Value = 0.3
$ws.Range("G47").Value = "lowperm"
$ws.Range("H47").Value = "20 md"

$ws.Range("A48").Value = "lowperm"
$ws.Range("B48").Value = "20 md"
$ws.Range("G48").Value = "leaked_perc"
$ws.Range("H48").Value = "0, 0.05, 0.1"

$ws.Range("A49").Value = "leaked_perc"
$ws.Range("B49").Value = "0, 0.05, 0.1"

# ---------------------------------------------------------------------------
# Restore the selection/scroll state shown in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("H38").Select()
